$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 2800.125
$ws.Cells.Item(33, 9).Value = 2400.2
$ws.Cells.Item(33, 10).Value = 3466.6667
$ws.Cells.Item(33, 11).Value = 2400.2
$ws.Cells.Item(33, 12).Value = 3466.6667
$ws.Cells.Item(33, 13).Value = -2171.2
$ws.Cells.Item(33, 14).Value = -3924.6667

$ws.Cells.Item(106, 8).Value = 150001500
$ws.Cells.Item(106, 9).Value = 200000670
$ws.Cells.Item(106, 10).Value = 4000
$ws.Cells.Item(106, 11).Value = 200000670
$ws.Cells.Item(106, 12).Value = 4000
$ws.Cells.Item(106, 13).Value = -200000039
$ws.Cells.Item(106, 14).Value = -5262

$ws.Cells.Item(138, 8).Value = 2033.33
$ws.Cells.Item(138, 10).Value = 2238.3562
$ws.Cells.Item(138, 12).Value = 6715.068600000001
$ws.Cells.Item(138, 14).Value = -16995.0686

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2560
$ws.Cells.Item(45, 9).Value = 2311.4546
$ws.Cells.Item(45, 11).Value = 2311.4546
$ws.Cells.Item(45, 13).Value = -1934.4546

$ws.Cells.Item(122, 8).Value = 93246.55
$ws.Cells.Item(122, 9).Value = 201522.4
$ws.Cells.Item(122, 10).Value = 3016.6667
$ws.Cells.Item(122, 11).Value = 604567.2
$ws.Cells.Item(122, 12).Value = 9050.000100000001
$ws.Cells.Item(122, 13).Value = -602117.2
$ws.Cells.Item(122, 14).Value = -13950.0001

$ws.Cells.Item(132, 8).Value = 3540.7368
$ws.Cells.Item(132, 9).Value = 3484.1052
$ws.Cells.Item(132, 10).Value = 3597.3684
$ws.Cells.Item(132, 11).Value = 10452.3156
$ws.Cells.Item(132, 12).Value = 10792.1052
$ws.Cells.Item(132, 13).Value = -7922.3156
$ws.Cells.Item(132, 14).Value = -15852.1052

$ws.Cells.Item(138, 8).Value = 140000
$ws.Cells.Item(138, 10).Value = 140000
$ws.Cells.Item(138, 12).Value = 140000
$ws.Cells.Item(138, 14).Value = -150280

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1900
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 13).Value = ""

$ws.Cells.Item(99, 8).Value = 966.6799999999999
$ws.Cells.Item(99, 9).Value = 923.1429000000001
$ws.Cells.Item(99, 10).Value = 1022.0909
$ws.Cells.Item(99, 11).Value = 923.1429000000001
$ws.Cells.Item(99, 12).Value = 1022.0909
$ws.Cells.Item(99, 13).Value = 574.8570999999999
$ws.Cells.Item(99, 14).Value = -4018.0909

$ws.Cells.Item(134, 8).Value = 2497
$ws.Cells.Item(134, 9).Value = 2406.9092
$ws.Cells.Item(134, 11).Value = 7220.7276
$ws.Cells.Item(134, 13).Value = -4685.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1837.8334
$ws.Cells.Item(58, 9).Value = 2127.4
$ws.Cells.Item(58, 10).Value = 1726.4615
$ws.Cells.Item(58, 11).Value = 2127.4
$ws.Cells.Item(58, 12).Value = 1726.4615
$ws.Cells.Item(58, 13).Value = -1924.4
$ws.Cells.Item(58, 14).Value = -2132.4615

$ws.Cells.Item(107, 8).Value = 2976966
$ws.Cells.Item(107, 9).Value = 5209048
$ws.Cells.Item(107, 10).Value = 856.6667
$ws.Cells.Item(107, 11).Value = 5209048
$ws.Cells.Item(107, 12).Value = 856.6667
$ws.Cells.Item(107, 13).Value = -5207128
$ws.Cells.Item(107, 14).Value = -4696.6667

$ws.Cells.Item(122, 8).Value = 1739.4166
$ws.Cells.Item(122, 9).Value = 1108
$ws.Cells.Item(122, 10).Value = 1999.4117
$ws.Cells.Item(122, 11).Value = 3324
$ws.Cells.Item(122, 12).Value = 5998.2351
$ws.Cells.Item(122, 13).Value = -874
$ws.Cells.Item(122, 14).Value = -10898.2351

$ws.Cells.Item(132, 8).Value = 10419657
$ws.Cells.Item(132, 9).Value = 2969.875
$ws.Cells.Item(132, 10).Value = 20836344
$ws.Cells.Item(132, 11).Value = 8909.625
$ws.Cells.Item(132, 12).Value = 62509032
$ws.Cells.Item(132, 13).Value = -6379.625
$ws.Cells.Item(132, 14).Value = -62514092

$ws.Cells.Item(136, 8).Value = 1837.8334
$ws.Cells.Item(136, 9).Value = 2127.4
$ws.Cells.Item(136, 10).Value = 1726.4615
$ws.Cells.Item(136, 11).Value = 6382.200000000001
$ws.Cells.Item(136, 12).Value = 5179.3845
$ws.Cells.Item(136, 13).Value = -3832.200000000001
$ws.Cells.Item(136, 14).Value = -10279.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(41, 8).Value = 2983.1667
$ws.Cells.Item(41, 10).Value = 3519.8
$ws.Cells.Item(41, 12).Value = 10559.4
$ws.Cells.Item(41, 14).Value = -11235.4

$ws.Cells.Item(54, 8).Value = 5081.0938
$ws.Cells.Item(54, 10).Value = 5081.0938
$ws.Cells.Item(54, 12).Value = 15243.2814
$ws.Cells.Item(54, 14).Value = -16361.2814

$ws.Cells.Item(125, 8).Value = 2138.3635
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 2138.3635
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 6415.0905
$ws.Cells.Item(125, 13).Value = ""
$ws.Cells.Item(125, 14).Value = -16255.0905

$ws.Cells.Item(131, 8).Value = 1011.14
$ws.Cells.Item(131, 10).Value = 1033.4791
$ws.Cells.Item(131, 12).Value = 3100.4373
$ws.Cells.Item(131, 14).Value = -13180.4373

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 111
$ws.Cells.Item(2, 9).Value = 50
$ws.Cells.Item(2, 10).Value = 172
$ws.Cells.Item(2, 11).Value = 50
$ws.Cells.Item(2, 12).Value = 172
$ws.Cells.Item(2, 13).Value = 63
$ws.Cells.Item(2, 14).Value = -398

$ws.Cells.Item(122, 8).Value = 3950.875
$ws.Cells.Item(122, 9).Value = 3003.5
$ws.Cells.Item(122, 10).Value = 4266.6665
$ws.Cells.Item(122, 11).Value = 9010.5
$ws.Cells.Item(122, 12).Value = 12799.9995
$ws.Cells.Item(122, 13).Value = -6560.5
$ws.Cells.Item(122, 14).Value = -17699.9995

$ws.Cells.Item(126, 8).Value = 1985
$ws.Cells.Item(126, 9).Value = 1979.3077
$ws.Cells.Item(126, 11).Value = 5937.9231
$ws.Cells.Item(126, 13).Value = -3467.9231

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 168800.67
$ws.Cells.Item(40, 9).Value = 334668
$ws.Cells.Item(40, 10).Value = 2933.3333
$ws.Cells.Item(40, 11).Value = 334668
$ws.Cells.Item(40, 12).Value = 2933.3333
$ws.Cells.Item(40, 13).Value = -334532
$ws.Cells.Item(40, 14).Value = -3205.3333

$ws.Cells.Item(68, 8).Value = 3237
$ws.Cells.Item(68, 9).Value = 2416.8572
$ws.Cells.Item(68, 10).Value = 4057.1428
$ws.Cells.Item(68, 11).Value = 2416.8572
$ws.Cells.Item(68, 12).Value = 4057.1428
$ws.Cells.Item(68, 13).Value = -1667.8572
$ws.Cells.Item(68, 14).Value = -5555.1428

$ws.Cells.Item(71, 8).Value = 3237
$ws.Cells.Item(71, 9).Value = 2416.8572
$ws.Cells.Item(71, 10).Value = 4057.1428
$ws.Cells.Item(71, 11).Value = 12084.286
$ws.Cells.Item(71, 12).Value = 20285.714
$ws.Cells.Item(71, 13).Value = -8340.286
$ws.Cells.Item(71, 14).Value = -27773.714

$ws.Cells.Item(122, 8).Value = 3235.3684
$ws.Cells.Item(122, 9).Value = 3134
$ws.Cells.Item(122, 11).Value = 9402
$ws.Cells.Item(122, 13).Value = -6952

$ws.Cells.Item(132, 8).Value = 2158.366
$ws.Cells.Item(132, 9).Value = 1245.8077
$ws.Cells.Item(132, 11).Value = 3737.4231
$ws.Cells.Item(132, 13).Value = -1207.4231

$ws.Cells.Item(136, 8).Value = 10419009
$ws.Cells.Item(136, 9).Value = 1982.6666
$ws.Cells.Item(136, 11).Value = 5947.9998
$ws.Cells.Item(136, 13).Value = -3397.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3044.3
$ws.Cells.Item(122, 9).Value = 1992.6
$ws.Cells.Item(122, 10).Value = 4096
$ws.Cells.Item(122, 11).Value = 5977.799999999999
$ws.Cells.Item(122, 12).Value = 12288
$ws.Cells.Item(122, 13).Value = -3527.799999999999
$ws.Cells.Item(122, 14).Value = -17188

$ws.Cells.Item(123, 8).Value = 24257.8
$ws.Cells.Item(123, 10).Value = 24257.8
$ws.Cells.Item(123, 12).Value = 24257.8
$ws.Cells.Item(123, 14).Value = -34057.8

$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).Value = ""

$ws.Cells.Item(135, 8).Value = 74690.836
$ws.Cells.Item(135, 10).Value = 74690.836
$ws.Cells.Item(135, 12).Value = 74690.836
$ws.Cells.Item(135, 14).Value = -84830.836
